# Fruta / hortaliza, semanal
# Insert one new weekly price record as a new row 400 in the Kiwi sheet,
# shifting the existing rows 400-425 down to 401-426.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 400 (Excel shifts rows 400:425 down to 401:426,
# and the new row 400 inherits formatting from the row above, matching the workbook's
# existing date-column style).
$ws.Rows(400).Insert()

# Populate the newly inserted row 400 with the new weekly record.
$ws.Range("A400").Value = 10
$ws.Range("B400").Value = "Vega Modelo de Temuco"
$ws.Range("C400").Value = "La Araucanía"
$ws.Range("D400").Value = 44706
$ws.Range("E400").Value = 9
$ws.Range("F400").Value = "Fruta"
$ws.Range("G400").Value = 100101
$ws.Range("H400").Value = "Berries"
$ws.Range("I400").Value = 100101007
$ws.Range("J400").Value = "Kiwi"
$ws.Range("K400").Value = "Hayward"
$ws.Range("L400").Value = "Especial"
$ws.Range("M400").Value = 115
$ws.Range("N400").Value = 17000
$ws.Range("O400").Value = 17000
$ws.Range("P400").Value = 17000
$ws.Range("Q400").Value = "$/bandeja 18 kilos"
$ws.Range("R400").Value = "Región de O'Higgins"
$ws.Range("S400").Value = 944
$ws.Range("T400").Value = 18
